$wb = $excel.ActiveWorkbook

# ALC row 18
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1785.7142
$ws.Range("I18").Value = 1785.7142
$ws.Range("K18").Value = 1785.7142
$ws.Range("M18").Value = -1501.7142

# ALC row 97
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H97").Value = 2000
$ws.Range("J97").Value = 2000
$ws.Range("L97").Value = 6000
$ws.Range("N97").Value = -6992

# ALC row 103
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H103").Value = 1897.6364
$ws.Range("I103").Value = 1999.5
$ws.Range("J103").Value = 1875
$ws.Range("K103").Value = 5998.5
$ws.Range("L103").Value = 5625
$ws.Range("M103").Value = -5412.5
$ws.Range("N103").Value = -6797

# ALC row 112
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 2359.818
$ws.Range("I112").Value = 1591.8
$ws.Range("J112").Value = 2999.8333
$ws.Range("K112").Value = 4775.4
$ws.Range("L112").Value = 8999.499899999999
$ws.Range("M112").Value = -3667.4
$ws.Range("N112").Value = -11215.4999

# ARM row 37
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 43919
$ws.Range("J37").Value = 43919
$ws.Range("L37").Value = 43919
$ws.Range("N37").Value = -44465

# ARM row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 43399
$ws.Range("I45").Value = 12331.667
$ws.Range("J45").Value = 90000
$ws.Range("K45").Value = 12331.667
$ws.Range("L45").Value = 90000
$ws.Range("M45").Value = -11954.667
$ws.Range("N45").Value = -90754

# ARM row 102
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 2111.4
$ws.Range("I102").Value = 2201.5
$ws.Range("K102").Value = 2201.5
$ws.Range("M102").Value = -579.5

# ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2909.2778
$ws.Range("I122").Value = 2528.4285
$ws.Range("J122").Value = 4242.25
$ws.Range("K122").Value = 7585.2855
$ws.Range("L122").Value = 12726.75
$ws.Range("M122").Value = -5135.2855
$ws.Range("N122").Value = -17626.75

# ARM row 129
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H129").Value = 47186.668
$ws.Range("J129").Value = 47186.668
$ws.Range("L129").Value = 47186.668
$ws.Range("N129").Value = -57186.668

# BSM row 20
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2504
$ws.Range("I20").Value = 2504
$ws.Range("K20").Value = 2504
$ws.Range("M20").Value = -2257

# BSM row 80
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 733.1111
$ws.Range("I80").Value = 51.25
$ws.Range("J80").Value = 1278.6
$ws.Range("K80").Value = 51.25
$ws.Range("L80").Value = 1278.6
$ws.Range("M80").Value = 946.75
$ws.Range("N80").Value = -3274.6

# BSM row 83
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H83").Value = 733.1111
$ws.Range("I83").Value = 51.25
$ws.Range("J83").Value = 1278.6
$ws.Range("K83").Value = 256.25
$ws.Range("L83").Value = 6393
$ws.Range("M83").Value = 4735.75
$ws.Range("N83").Value = -16377

# BSM row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 5649.3
$ws.Range("J86").Value = 8373.75
$ws.Range("L86").Value = 8373.75
$ws.Range("N86").Value = -10619.75

# BSM row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 5649.3
$ws.Range("J89").Value = 8373.75
$ws.Range("L89").Value = 41868.75
$ws.Range("N89").Value = -53100.75

# BSM row 94
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 4322
$ws.Range("I94").Value = 2583.6
$ws.Range("K94").Value = 2583.6
$ws.Range("M94").Value = -2132.6

# BSM row 105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 17501.75
$ws.Range("I105").Value = 17501.75
$ws.Range("K105").Value = 17501.75
$ws.Range("M105").Value = -15754.75

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2989.75
$ws.Range("I134").Value = 2989.75
$ws.Range("K134").Value = 8969.25
$ws.Range("M134").Value = -6434.25

# CRP row 10
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 525.4
$ws.Range("J10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("N10").ClearContents()

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2659.3333
$ws.Range("I31").Value = 2239.3333
$ws.Range("J31").Value = 3499.3333
$ws.Range("K31").Value = 2239.3333
$ws.Range("L31").Value = 3499.3333
$ws.Range("M31").Value = -1944.3333
$ws.Range("N31").Value = -4089.3333

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2659.3333
$ws.Range("I34").Value = 2239.3333
$ws.Range("J34").Value = 3499.3333
$ws.Range("K34").Value = 2239.3333
$ws.Range("L34").Value = 3499.3333
$ws.Range("M34").Value = -2037.3333
$ws.Range("N34").Value = -3903.3333

# CRP row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3621.875
$ws.Range("I58").Value = 3865
$ws.Range("K58").Value = 3865
$ws.Range("M58").Value = -3662

# CRP row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 3621.875
$ws.Range("I136").Value = 3865
$ws.Range("K136").Value = 11595
$ws.Range("M136").Value = -9045

# CUL row 107
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 585.7143

# CUL row 132
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 3375
$ws.Range("I132").Value = 3000
$ws.Range("K132").Value = 27000
$ws.Range("M132").Value = -24470

# GSM row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 8179.8
$ws.Range("I80").Value = 3450
$ws.Range("J80").Value = 11333
$ws.Range("K80").Value = 3450
$ws.Range("L80").Value = 11333
$ws.Range("M80").Value = -2452
$ws.Range("N80").Value = -13329

# GSM row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 8179.8
$ws.Range("I83").Value = 3450
$ws.Range("J83").Value = 11333
$ws.Range("K83").Value = 17250
$ws.Range("L83").Value = 56665
$ws.Range("M83").Value = -12258
$ws.Range("N83").Value = -66649

# GSM row 97
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2552
$ws.Range("I97").Value = 1613.3334
$ws.Range("K97").Value = 1613.3334
$ws.Range("M97").Value = -1117.3334

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3496.75
$ws.Range("I132").Value = 2307.875
$ws.Range("J132").Value = 5874.5
$ws.Range("K132").Value = 6923.625
$ws.Range("L132").Value = 17623.5
$ws.Range("M132").Value = -4393.625
$ws.Range("N132").Value = -22683.5

# LTW row 16
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1806.1666
$ws.Range("I16").Value = 967.6
$ws.Range("K16").Value = 967.6
$ws.Range("M16").Value = -797.6

# LTW row 19
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("M19").ClearContents()

# LTW row 55
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 734.7857
$ws.Range("I55").Value = 787.36365
$ws.Range("J55").Value = 542
$ws.Range("K55").Value = 787.36365
$ws.Range("L55").Value = 542
$ws.Range("M55").Value = -614.36365
$ws.Range("N55").Value = -888

# LTW row 82
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2700.7273
$ws.Range("J82").Value = 5501.3335
$ws.Range("L82").Value = 5501.3335
$ws.Range("N82").Value = -6223.3335

# LTW row 85
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 2700.7273
$ws.Range("J85").Value = 5501.3335
$ws.Range("L85").Value = 5501.3335
$ws.Range("N85").Value = -7997.3335

# WVR row 81
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3597
$ws.Range("I81").Value = 1435.8
$ws.Range("J81").Value = 9000
$ws.Range("K81").Value = 2871.6
$ws.Range("L81").Value = 18000
$ws.Range("M81").Value = -1810.6
$ws.Range("N81").Value = -20122

# WVR row 84
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 3597
$ws.Range("I84").Value = 1435.8
$ws.Range("J84").Value = 9000
$ws.Range("K84").Value = 14358
$ws.Range("L84").Value = 90000
$ws.Range("M84").Value = -9054
$ws.Range("N84").Value = -100608

# WVR row 96
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2136.182
$ws.Range("I96").Value = 1989.8
$ws.Range("K96").Value = 1989.8
$ws.Range("M96").Value = -616.8
